# "Generate Report for Archive"
#
# 1. The shared status string "Ready for handoff" becomes "In Translation".
#    This string is shared by the Overview sheet (columns zh-cn/de-de, i.e.
#    E2:E3 and F2:F3) as well as the per-locale "Status" column (column C,
#    rows 2-3) on both the zh-cn and de-de sheets. Using Replace() updates
#    every occurrence consistently (and in-place as a single shared string)
#    instead of touching each cell independently.
# 2. Because the status text got shorter, the Status/locale columns were
#    re-sized (narrower) on all three sheets.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
